# Applies the "Atualizado por script em 31-10-2023 15:01" update to the
# Lebanon Premier League 2023-2024 results sheet:
#   1) Rows 16 and 17 (match #15 Tripoli-AlGhazieh and match #16 AlSahel-AlAnsar)
#      had their match order swapped upstream, so the two rows exchange all of
#      their match-specific data (home/away teams, odds, timestamps, url).
#   2) Five new match rows (35-39, Indice 34-38) are appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($ws, [int]$row, $values)
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}

# ---------------------------------------------------------------------------
# 1) Swap the contents of rows 16 and 17 for the columns that differ between
#    the two matches (F,H,J,L,M,N,P,Q,R,T,U,V). Columns A-E,G,I,K,O,S already
#    hold identical/unaffected values for both rows.
# ---------------------------------------------------------------------------

$row16New = @{
    F = "Al Sahel"
    H = "Al Ansar"
    J = 3.98
    L = 4.23
    M = "20/08/2023 14:06"
    N = 3.39
    P = 3.49
    Q = "20/08/2023 15:31"
    R = 1.76
    T = 1.82
    U = "20/08/2023 14:06"
    V = "https://www.betexplorer.com/football/lebanon/premier-league/al-sahel-al-ansar/SMIGhbSr/"
}

$row17New = @{
    F = "Tripoli"
    H = "Al Ghazieh"
    J = 2
    L = 1.67
    M = "20/08/2023 15:00"
    N = 3.08
    P = 3.64
    Q = "20/08/2023 15:00"
    R = 3.46
    T = 5.06
    U = "20/08/2023 15:00"
    V = "https://www.betexplorer.com/football/lebanon/premier-league/tripoli-sc-al-ghazieh/QDWlcvKR/"
}

Set-RowValues $ws 16 $row16New
Set-RowValues $ws 17 $row17New

# ---------------------------------------------------------------------------
# 2) Append the five new match rows (35-39) after the existing last row (34).
#    Formatting (bold/border/center on column A, date number format on
#    column E) is copied down from row 34 before the values are written.
# ---------------------------------------------------------------------------

$newRows = @(
    @{
        A = 34; B = "lebanon"; C = "premier-league"; D = "2023-2024"
        E = 45226.58333333334
        F = "Al Sahel"; G = 1; H = "Bourj FC"; I = 2
        J = 2.32; K = "26/10/2023 02:12"
        L = 3.05; M = "27/10/2023 13:04"
        N = 3;    O = "26/10/2023 02:12"
        P = 2.68; Q = "27/10/2023 13:04"
        R = 2.91; S = "26/10/2023 02:12"
        T = 2.71; U = "27/10/2023 13:04"
        V = "https://www.betexplorer.com/football/lebanon/premier-league/al-sahel-bourj/xz7yJftl/"
    },
    @{
        A = 35; B = "lebanon"; C = "premier-league"; D = "2023-2024"
        E = 45227.58333333334
        F = "Al Ghazieh"; G = 1; H = "Al Hikma"; I = 0
        J = 2.59; K = "27/10/2023 02:13"
        L = 3.02; M = "28/10/2023 13:23"
        N = 3.09; O = "27/10/2023 02:13"
        P = 3.16; Q = "28/10/2023 13:23"
        R = 2.47; S = "27/10/2023 02:13"
        T = 2.36; U = "28/10/2023 13:23"
        V = "https://www.betexplorer.com/football/lebanon/premier-league/al-ghazieh-al-hikma/jwBuIEef/"
    },
    @{
        A = 36; B = "lebanon"; C = "premier-league"; D = "2023-2024"
        E = 45227.625
        F = "Al Ansar"; G = 1; H = "Tadamon"; I = 0
        J = 1.31; K = "27/10/2023 03:12"
        L = 1.36; M = "27/10/2023 22:05"
        N = 4.53; O = "27/10/2023 03:12"
        P = 4.49; Q = "28/10/2023 13:03"
        R = 7.78; S = "27/10/2023 03:12"
        T = 8.56; U = "27/10/2023 22:05"
        V = "https://www.betexplorer.com/football/lebanon/premier-league/al-ansar-tadamon/p0AqHYA0/"
    },
    @{
        A = 37; B = "lebanon"; C = "premier-league"; D = "2023-2024"
        E = 45228.55208333334
        F = "Racing"; G = 3; H = "Al Ahli Nabatiya"; I = 0
        J = 1.75; K = "28/10/2023 02:42"
        L = 1.72; M = "29/10/2023 11:56"
        N = 3.44; O = "28/10/2023 02:42"
        P = 3.56; Q = "29/10/2023 11:56"
        R = 3.96; S = "28/10/2023 02:42"
        T = 4.79; U = "29/10/2023 11:56"
        V = "https://www.betexplorer.com/football/lebanon/premier-league/racing-al-ahli-nabatiya/42QpFCuD/"
    },
    @{
        A = 38; B = "lebanon"; C = "premier-league"; D = "2023-2024"
        E = 45228.65625
        F = "Safa"; G = 1; H = "Al Ahed"; I = 3
        J = 6.22; K = "28/10/2023 05:12"
        L = 4.95; M = "29/10/2023 10:47"
        N = 4.49; O = "28/10/2023 05:12"
        P = 4.53; Q = "29/10/2023 13:47"
        R = 1.36; S = "28/10/2023 05:12"
        T = 1.53; U = "29/10/2023 10:47"
        V = "https://www.betexplorer.com/football/lebanon/premier-league/safa-al-ahed/rgPlEWfJ/"
    }
)

$targetRow = 35
foreach ($rowData in $newRows) {
    # Copy formatting from the last existing data row so the new row keeps
    # the same cell styles (bold/border/center on A, date format on E).
    $ws.Range("A34:V34").Copy()
    $ws.Range("A${targetRow}:V${targetRow}").PasteSpecial(-4122)

    Set-RowValues $ws $targetRow $rowData
    $targetRow++
}

$excel.CutCopyMode = 0
